$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update category labels ---
$ws.Range("J4").Value = "n=100"
$ws.Range("J5").Value = "n=1000"
$ws.Range("J6").Value = "n=2000"

# --- Update existing series headers (unchanged text, re-assert for safety) ---
$ws.Range("K3").Value = "bubble sort"
$ws.Range("L3").Value = "Insert sort"
$ws.Range("M3").Value = "Select sort"
$ws.Range("N3").Value = "Cocktail Sort"

# --- Update existing series data (bubble, insert, select, cocktail) ---
$ws.Range("K4").Value = 0.019
$ws.Range("K5").Value = 0.231
$ws.Range("K6").Value = 0.418

$ws.Range("L4").Value = 0.023
$ws.Range("L5").Value = 0.275
$ws.Range("L6").Value = 0.454

$ws.Range("M4").Value = 0.02
$ws.Range("M5").Value = 0.3
$ws.Range("M6").Value = 0.5

$ws.Range("N4").Value = 0.02
$ws.Range("N5").Value = 0.244
$ws.Range("N6").Value = 0.6

# --- Add the 3 new algorithm columns ---
$ws.Range("O3").Value = "Merge Sort"
$ws.Range("O4").Value = 0.014
$ws.Range("O5").Value = 0.222
$ws.Range("O6").Value = 0.404

$ws.Range("P3").Value = "Quick Sort"
$ws.Range("P4").Value = 0.019
$ws.Range("P5").Value = 0.205
$ws.Range("P6").Value = 0.427

$ws.Range("Q3").Value = "Heap Sort"
$ws.Range("Q4").Value = 0.022
$ws.Range("Q5").Value = 0.257
$ws.Range("Q6").Value = 0.491

# --- Update the sheet view: zoom + selection ---
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("K24").Select()

# --- Chart: add the 3 new series (Merge Sort, Quick Sort, Heap Sort) ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$merge = $chart.SeriesCollection().NewSeries()
$merge.Name = "=Hoja1!`$O`$3"
$merge.XValues = "=Hoja1!`$J`$4:`$J`$6"
$merge.Values = "=Hoja1!`$O`$4:`$O`$6"
$merge.Smooth = $false
$merge.Border.Weight = 2.75
$merge.Border.Color = 0x7399

$quick = $chart.SeriesCollection().NewSeries()
$quick.Name = "=Hoja1!`$P`$3"
$quick.XValues = "=Hoja1!`$J`$4:`$J`$6"
$quick.Values = "=Hoja1!`$P`$4:`$P`$6"
$quick.Smooth = $false
$quick.Border.Weight = 2.75
$quick.Border.Color = 0x2B6843

$heap = $chart.SeriesCollection().NewSeries()
$heap.Name = "=Hoja1!`$Q`$3"
$heap.XValues = "=Hoja1!`$J`$4:`$J`$6"
$heap.Values = "=Hoja1!`$Q`$4:`$Q`$6"
$heap.Smooth = $false
$heap.Border.Weight = 2.75
$heap.Border.Color = 0x5A97F1

# --- Resize/reposition the chart to its new extent ---
$co.Top = 14.1
$co.Left = 59.6375

Write-Host "Edit complete"
